# Document upload functionality: add an INSERT-statement helper column (C)
# that builds a SQL upsert for every building-code row, and clean up a
# handful of description strings that used parentheses inside the text
# (replaced with angle brackets so they don't collide with the SQL quoting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix up the descriptions that previously used a parenthetical aside
#    (some of them were even split across two lines). Replace with a
#    single-line version using angle brackets instead of parentheses.
# ---------------------------------------------------------------------
$ws.Range("B23").Value = "EWB-10-Conditioning Building <Small Scale Propellant Assembly>"
$ws.Range("B28").Value = "GFB-02-Main Entrance Building <pedestrian>"
$ws.Range("B30").Value = "GFB-04-Cafeteria <Dining Hall Building>"
$ws.Range("B37").Value = "GFB-11-Security Tower <8 units>"
$ws.Range("B48").Value = "GFB-14-Water Tank Building <Potable water>"
$ws.Range("B55").Value = "ISB-02-Inert Material Storage Building <Incoming Inspection and Calibration>"

# These two got long enough (after losing the line break) that they need
# to wrap within the existing column width.
$ws.Range("B23").WrapText = $true
$ws.Range("B55").WrapText = $true
$ws.Range("B55").RowHeight = 28.5

# ---------------------------------------------------------------------
# 2. Build column C: one SQL INSERT statement per row, generated from the
#    code/description columns. Enter the formula in C2 first, then fill
#    it down through C64 (mirrors typing it once and dragging the fill
#    handle down, which is what produces a shared formula group C3:C64).
# ---------------------------------------------------------------------
$formulaC2 = "=""INSERT INTO building_codes (code, description) VALUES ( '""&A2&""', '""&B2&""')""&"" ON CONFLICT(code) DO NOTHING;"""
$ws.Range("C2").Formula = $formulaC2

$formulaC3 = "=""INSERT INTO building_codes (code, description) VALUES ( '""&A3&""', '""&B3&""')""&"" ON CONFLICT(code) DO NOTHING;"""
$ws.Range("C3:C64").Formula = $formulaC3

# ---------------------------------------------------------------------
# 3. Leave the sheet selected on the new helper column, like it was left
#    after building it.
# ---------------------------------------------------------------------
$ws.Range("C2:C64").Select()
